# Weekly update: insert a new Zanahoria price record (Terminal Hortofrutícola
# Agro Chillán) above row 156, shifting the existing historical rows down by
# one, and populate the newly-inserted row with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 156; this pushes the old rows 156..211
# down to 157..212 (carrying their existing values/formatting with them).
$ws.Rows.Item(156).Insert()

# Populate the newly inserted row 156 with this week's data.
$ws.Cells.Item(156, 1).Value  = 7
$ws.Cells.Item(156, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(156, 3).Value  = "Ñuble"
$ws.Cells.Item(156, 4).Value  = 44524
$ws.Cells.Item(156, 5).Value  = 16
$ws.Cells.Item(156, 6).Value  = 100114013
$ws.Cells.Item(156, 7).Value  = "Zanahoria"
$ws.Cells.Item(156, 8).Value  = "Sin especificar"
$ws.Cells.Item(156, 9).Value  = "Primera"
$ws.Cells.Item(156, 10).Value = 60
$ws.Cells.Item(156, 11).Value = 8000
$ws.Cells.Item(156, 12).Value = 8500
$ws.Cells.Item(156, 13).Value = 8250
$ws.Cells.Item(156, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(156, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(156, 16).Value = 412
$ws.Cells.Item(156, 17).Value = 20
$ws.Cells.Item(156, 18).Value = "Hortaliza"
